$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.1655629139072848
$ws.Range("C2").Value = 0.6059602649006622
$ws.Range("J2").Value = 0.01655629139072848
$ws.Range("P2").Value = 0.1423841059602649
$ws.Range("S2").Value = 0.0695364238410596

# Row 3
$ws.Range("B3").Value = 0.01052631578947368
$ws.Range("C3").Value = 0.03157894736842105
$ws.Range("J3").Value = 0.02631578947368421
$ws.Range("P3").Value = 0.7315789473684211
$ws.Range("S3").Value = 0.2

# Row 4
$ws.Range("J4").Value = 0.08695652173913043
$ws.Range("P4").Value = 0.6521739130434783
$ws.Range("S4").Value = 0.2608695652173913

# Row 6
$ws.Range("B6").Value = 0.0867579908675799
$ws.Range("D6").Value = 0.0273972602739726
$ws.Range("F6").Value = 0.0593607305936073
$ws.Range("J6").Value = 0.2100456621004566
$ws.Range("O6").Value = 0.0319634703196347
$ws.Range("Q6").Value = 0.1598173515981735
$ws.Range("R6").Value = 0.0684931506849315
$ws.Range("S6").Value = 0.3561643835616438

# Row 7
$ws.Range("B7").Value = 0.09836065573770492
$ws.Range("D7").Value = 0.02185792349726776
$ws.Range("F7").Value = 0.04918032786885246
$ws.Range("J7").Value = 0.1202185792349727
$ws.Range("O7").Value = 0.01092896174863388
$ws.Range("Q7").Value = 0.2404371584699453
$ws.Range("R7").Value = 0.1038251366120219
$ws.Range("S7").Value = 0.3551912568306011

# Row 8
$ws.Range("B8").Value = 0.1159737417943107
$ws.Range("D8").Value = 0.01312910284463895
$ws.Range("F8").Value = 0.0700218818380744
$ws.Range("J8").Value = 0.1094091903719912
$ws.Range("O8").Value = 0.03282275711159737
$ws.Range("Q8").Value = 0.1838074398249453
$ws.Range("R8").Value = 0.07221006564551423
$ws.Range("S8").Value = 0.4026258205689278

# Row 9
$ws.Range("B9").Value = 0.09090909090909091
$ws.Range("D9").Value = 0.0160427807486631
$ws.Range("F9").Value = 0.1122994652406417
$ws.Range("J9").Value = 0.08021390374331551
$ws.Range("O9").Value = 0.0160427807486631
$ws.Range("Q9").Value = 0.1925133689839572
$ws.Range("R9").Value = 0.06951871657754011
$ws.Range("S9").Value = 0.4224598930481284

# Row 10
$ws.Range("B10").Value = 0.1203389830508475
$ws.Range("D10").Value = 0.02457627118644068
$ws.Range("E10").Value = 0.000847457627118644
$ws.Range("F10").Value = 0.06440677966101695
$ws.Range("J10").Value = 0.1059322033898305
$ws.Range("O10").Value = 0.01610169491525424
$ws.Range("Q10").Value = 0.2177966101694915
$ws.Range("R10").Value = 0.08135593220338982
$ws.Range("S10").Value = 0.3686440677966102

# Row 11
$ws.Range("G11").Value = 0.1453900709219858
$ws.Range("J11").Value = 0.124113475177305
$ws.Range("K11").Value = 0.2234042553191489
$ws.Range("L11").Value = 0.4964539007092199
$ws.Range("S11").Value = 0.01063829787234043

# Row 12
$ws.Range("G12").Value = 0.823943661971831
$ws.Range("J12").Value = 0.1197183098591549
$ws.Range("K12").Value = 0.01408450704225352
$ws.Range("L12").Value = 0.02112676056338028
$ws.Range("S12").Value = 0.02112676056338028

# Row 13
$ws.Range("G13").Value = 0.6808510638297872
$ws.Range("J13").Value = 0.2340425531914894
$ws.Range("S13").Value = 0.0851063829787234

# Row 15
$ws.Range("F15").Value = 0.0187793427230047
$ws.Range("H15").Value = 0.1784037558685446
$ws.Range("I15").Value = 0.07511737089201878
$ws.Range("J15").Value = 0.3380281690140845
$ws.Range("K15").Value = 0.06103286384976526
$ws.Range("M15").Value = 0.02816901408450704
$ws.Range("O15").Value = 0.05633802816901409
$ws.Range("S15").Value = 0.244131455399061

# Row 16
$ws.Range("F16").Value = 0.01951219512195122
$ws.Range("H16").Value = 0.1902439024390244
$ws.Range("I16").Value = 0.05853658536585366
$ws.Range("J16").Value = 0.375609756097561
$ws.Range("K16").Value = 0.0975609756097561
$ws.Range("M16").Value = 0.02439024390243903
$ws.Range("O16").Value = 0.05853658536585366
$ws.Range("S16").Value = 0.175609756097561

# Row 17
$ws.Range("F17").Value = 0.01978021978021978
$ws.Range("H17").Value = 0.1692307692307692
$ws.Range("I17").Value = 0.1032967032967033
$ws.Range("J17").Value = 0.4417582417582417
$ws.Range("K17").Value = 0.0945054945054945
$ws.Range("M17").Value = 0.02637362637362637
$ws.Range("N17").Value = 0.002197802197802198
$ws.Range("O17").Value = 0.04615384615384616
$ws.Range("S17").Value = 0.09670329670329671

# Row 18
$ws.Range("F18").Value = 0.02272727272727273
$ws.Range("H18").Value = 0.1761363636363636
$ws.Range("I18").Value = 0.09659090909090909
$ws.Range("J18").Value = 0.4431818181818182
$ws.Range("K18").Value = 0.07386363636363637
$ws.Range("M18").Value = 0.02272727272727273
$ws.Range("O18").Value = 0.02840909090909091
$ws.Range("S18").Value = 0.1363636363636364

# Row 19
$ws.Range("F19").Value = 0.01231527093596059
$ws.Range("H19").Value = 0.2241379310344828
$ws.Range("I19").Value = 0.08045977011494253
$ws.Range("J19").Value = 0.3538587848932677
$ws.Range("K19").Value = 0.1026272577996716
$ws.Range("M19").Value = 0.0180623973727422
$ws.Range("O19").Value = 0.07881773399014778
$ws.Range("S19").Value = 0.1297208538587849
